{"js": "// Remove \"blastula, keyring, \" from the list of R packages that need to be\n// manually installed (the gmail-notification functions that relied on\n// blastula/keyring were removed from the project).\nconst body = context.document.body;\n\nconst searchResults = body.search(\", blastula, keyring, roxygen2, and \", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\n\nsearchResults.items[0].insertText(\", roxygen2, and \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Remove \"blastula, keyring, \" from the list of R packages that need to be\n# manually installed (the gmail-notification functions that relied on\n# blastula/keyring were removed from the project).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \", blastula, keyring, roxygen2, and \"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \", roxygen2, and \"\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute(\n    $find.Text,\n    $find.MatchCase,\n    $find.MatchWholeWord,\n    $find.MatchWildcards,\n    $null, $null, $find.Forward, $find.Wrap, $null,\n    $find.Replacement.Text,\n    2  # wdReplaceAll\n) | Out-Null\n"}
